$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.951.16"
$ws.Range("E2").Value = "  +4.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.81"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.76"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.98"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.422"
$ws.Range("E9").Value = "  +4.54%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0939"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.62"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.609.91"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("E14").Value = "  +8.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.56"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.80"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.808"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.282.52"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.818.73"
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.08"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.18"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  +6.18%  "
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.21"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.50"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.05"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0657"
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.45"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.74"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000225"
$ws.Range("E43").Value = "  -16.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0980"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.22"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.15"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.41"
$ws.Range("E47").Value = "  -8.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.477.08"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.61"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.73"
$ws.Range("E51").Value = "  -2.65%  "
